# Update recalculated NATMI edge-weight values for the Mif-Cxcr4 LR pair
# after refreshing the TPM inputs (see commit "update scripts wuth new tpm").
# Only the numeric value cells changed; row/column layout is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value = 6.556445
$ws.Range("H2").Value = 19.669335
$ws.Range("I2").Value = 0.2003595613103873
$ws.Range("J2").Value = 0.2003595613103873
$ws.Range("M2").Value = 13.71977066666667
$ws.Range("N2").Value = 41.159312
$ws.Range("O2").Value = 0.5515038136402627
$ws.Range("P2").Value = 0.5515038136402626
$ws.Range("Q2").Value = 89.95292178861334
$ws.Range("R2").Value = 809.57629609752
$ws.Range("S2").Value = 0.1104990621619686
$ws.Range("T2").Value = 0.1104990621619686

# Row 3 (ECs -> FAPs)
$ws.Range("G3").Value = 6.556445
$ws.Range("H3").Value = 19.669335
$ws.Range("I3").Value = 0.2003595613103873
$ws.Range("J3").Value = 0.2003595613103873
$ws.Range("O3").Value = 0.172077867958883
$ws.Range("P3").Value = 0.1720778679588829
$ws.Range("Q3").Value = 28.06672703836167
$ws.Range("R3").Value = 252.600543345255
$ws.Range("S3").Value = 0.03447744613546853
$ws.Range("T3").Value = 0.03447744613546853

# Row 4 (ECs -> MuSCs)
$ws.Range("G4").Value = 6.556445
$ws.Range("H4").Value = 19.669335
$ws.Range("I4").Value = 0.2003595613103873
$ws.Range("J4").Value = 0.2003595613103873
$ws.Range("O4").Value = 0.2764183184008545
$ws.Range("P4").Value = 0.2764183184008545
$ws.Range("Q4").Value = 45.08515582499833
$ws.Range("R4").Value = 405.766402424985
$ws.Range("S4").Value = 0.05538305301295015
$ws.Range("T4").Value = 0.05538305301295015

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.3842320902647997
$ws.Range("J5").Value = 0.3842320902647997
$ws.Range("M5").Value = 13.71977066666667
$ws.Range("N5").Value = 41.159312
$ws.Range("O5").Value = 0.5515038136402627
$ws.Range("P5").Value = 0.5515038136402626
$ws.Range("Q5").Value = 172.5038672385689
$ws.Range("R5").Value = 1552.53480514712
$ws.Range("S5").Value = 0.2119054631040067
$ws.Range("T5").Value = 0.2119054631040067

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.3842320902647997
$ws.Range("J6").Value = 0.3842320902647997
$ws.Range("O6").Value = 0.172077867958883
$ws.Range("P6").Value = 0.1720778679588829
$ws.Range("S6").Value = 0.06611783889415179
$ws.Range("T6").Value = 0.06611783889415179

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.3842320902647997
$ws.Range("J7").Value = 0.3842320902647997
$ws.Range("O7").Value = 0.2764183184008545
$ws.Range("P7").Value = 0.2764183184008545
$ws.Range("S7").Value = 0.1062087882666412
$ws.Range("T7").Value = 0.1062087882666412

# Row 8 (MuSCs -> ECs)
$ws.Range("G8").Value = 13.59357133333334
$ws.Range("I8").Value = 0.4154083484248129
$ws.Range("J8").Value = 0.415408348424813
$ws.Range("M8").Value = 13.71977066666667
$ws.Range("N8").Value = 41.159312
$ws.Range("O8").Value = 0.5515038136402627
$ws.Range("P8").Value = 0.5515038136402626
$ws.Range("Q8").Value = 186.5006812343076
$ws.Range("R8").Value = 1678.506131108768
$ws.Range("S8").Value = 0.2290992883742873
$ws.Range("T8").Value = 0.2290992883742873

# Row 9 (MuSCs -> FAPs)
$ws.Range("G9").Value = 13.59357133333334
$ws.Range("I9").Value = 0.4154083484248129
$ws.Range("J9").Value = 0.415408348424813
$ws.Range("O9").Value = 0.172077867958883
$ws.Range("P9").Value = 0.1720778679588829
$ws.Range("Q9").Value = 58.19114719778245
$ws.Range("R9").Value = 523.720324780042
$ws.Range("S9").Value = 0.07148258292926261
$ws.Range("T9").Value = 0.07148258292926261

# Row 10 (MuSCs -> MuSCs)
$ws.Range("G10").Value = 13.59357133333334
$ws.Range("I10").Value = 0.4154083484248129
$ws.Range("J10").Value = 0.415408348424813
$ws.Range("O10").Value = 0.2764183184008545
$ws.Range("P10").Value = 0.2764183184008545
$ws.Range("Q10").Value = 93.47569937390823
$ws.Range("R10").Value = 841.2812943651739
$ws.Range("S10").Value = 0.114826477121263
$ws.Range("T10").Value = 0.1148264771212631
